$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.663.43"
$ws.Range("E2").Value = "  +3.15%  "
$ws.Range("D3").Value = "2.189.01"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'257.93"
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("D6").Value = "'80.71"
$ws.Range("E6").Value = "  +9.47%  "
$ws.Range("E7").Value = "  +2.62%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.591"
$ws.Range("E9").Value = "  +1.57%  "
$ws.Range("D10").Value = "'42.64"
$ws.Range("E10").Value = "  +5.18%  "
$ws.Range("D11").Value = "'0.0917"
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("E12").Value = "  +2.49%  "
$ws.Range("E13").Value = "  +2.22%  "
$ws.Range("D14").Value = "2.516.38"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").Value = "'14.20"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "2.188.70"
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("D17").Value = "'0.774"
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("D18").Value = "43.566.49"
$ws.Range("E18").Value = "  +3.11%  "
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").Value = "'70.03"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").Value = "'5.91"
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("D22").Value = "'2.38"
$ws.Range("E22").Value = "  +11.26%  "
$ws.Range("D23").Value = "'229.68"
$ws.Range("E23").Value = "  +1.29%  "
$ws.Range("D24").Value = "'8.92"
$ws.Range("E24").Value = "  -5.43%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "'41.20"
$ws.Range("E26").Value = "  +13.19%  "
$ws.Range("D27").Value = "'10.59"
$ws.Range("E27").Value = "  +1.57%  "
$ws.Range("D28").Value = "'3.36"
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'2.21"
$ws.Range("E29").Value = "  +2.00%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.18"
$ws.Range("E30").Value = "  -1.64%  "
$ws.Range("D31").Value = "'172.60"
$ws.Range("E31").Value = "  +1.54%  "
$ws.Range("D32").Value = "'20.32"
$ws.Range("E32").Value = "  +1.73%  "
$ws.Range("D33").Value = "'0.0870"
$ws.Range("E33").Value = "  +7.56%  "
$ws.Range("D34").Value = "'5.24"
$ws.Range("E34").Value = "  +2.57%  "
$ws.Range("E35").Value = "  +5.95%  "
$ws.Range("D36").Value = "'0.122"
$ws.Range("E36").Value = "  +1.37%  "
$ws.Range("D37").Value = "'4.46"
$ws.Range("E37").Value = "  +5.57%  "
$ws.Range("D38").Value = "'0.0353"
$ws.Range("E38").Value = "  +5.18%  "
$ws.Range("D39").Value = "'13.06"
$ws.Range("E39").Value = "  +11.05%  "
$ws.Range("D40").Value = "'2.83"
$ws.Range("E40").Value = "  +16.80%  "
$ws.Range("D41").Value = "'2.09"
$ws.Range("E41").Value = "  +1.85%  "
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").Value = "'62.10"
$ws.Range("E42").Value = "  +4.66%  "
$ws.Range("B43").Value = "THORChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D43").Value = "'5.42"
$ws.Range("E43").Value = "  +5.62%  "
$ws.Range("E44").Value = "  +0.78%  "
$ws.Range("D45").Value = "'100.67"
$ws.Range("E45").Value = "  -2.03%  "
$ws.Range("D46").Value = "'0.0983"
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").Value = "'1.17"
$ws.Range("E48").Value = "  +3.87%  "
$ws.Range("E49").Value = "  +26.97%  "
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").Value = "'1.10"
$ws.Range("E50").Value = "  +1.61%  "
$ws.Range("D51").Value = "'0.436"
$ws.Range("E51").Value = "  -7.53%  "
